$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.309.82"
$ws.Range("E2").Value = "  +1.45%  "

# Row 3
$ws.Range("D3").Value = "3.781.27"
$ws.Range("E3").Value = "  -0.47%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "630.80"
$ws.Range("E5").Value = "  +3.83%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.38"
$ws.Range("E6").Value = "  +1.74%  "

# Row 7
$ws.Range("D7").Value = "3.779.44"
$ws.Range("E7").Value = "  -0.41%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  +1.10%  "

# Row 10
$ws.Range("E10").Value = "  -0.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.461"
$ws.Range("E11").Value = "  +2.46%  "

# Row 12
$ws.Range("E12").Value = "  -2.37%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("E13").Value = "  -1.56%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.93"
$ws.Range("E14").Value = "  -0.51%  "

# Row 15
$ws.Range("D15").Value = "4.414.14"
$ws.Range("E15").Value = "  -0.40%  "

# Row 16
$ws.Range("D16").Value = "3.778.56"
$ws.Range("E16").Value = "  -0.25%  "

# Row 17
$ws.Range("D17").Value = "69.300.08"
$ws.Range("E17").Value = "  +1.49%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.64"
$ws.Range("E18").Value = "  -2.51%  "

# Row 19
$ws.Range("E19").Value = "  +0.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.02"
$ws.Range("E20").Value = "  -0.55%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "464.65"
$ws.Range("E21").Value = "  +0.45%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.55"
$ws.Range("E22").Value = "  -0.53%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.709"
$ws.Range("E23").Value = "  +1.31%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.06"
$ws.Range("E24").Value = "  -0.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("E25").Value = "  -1.56%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.05"
$ws.Range("E26").Value = "  +0.46%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.15"
$ws.Range("E27").Value = "  +1.96%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("E28").Value = "  +0.77%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("D30").Value = "3.929.89"
$ws.Range("E30").Value = "  -0.43%  "

# Row 31
$ws.Range("E31").Value = "  +2.74%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.68"
$ws.Range("E32").Value = "  +1.74%  "

# Row 33
$ws.Range("E33").Value = "  -1.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.61"
$ws.Range("E34").Value = "  -1.69%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.171"
$ws.Range("E35").Value = "  +14.78%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.18%  "

# Row 37
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.01"
$ws.Range("E37").Value = "  -0.30%  "

# Row 38
$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "3.733.61"
$ws.Range("E38").Value = "  -0.25%  "

# Row 39
$ws.Range("E39").Value = "  +0.22%  "

# Row 40
$ws.Range("E40").Value = "  +3.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.81"
$ws.Range("E41").Value = "  -1.11%  "

# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.961"
$ws.Range("E42").Value = "  -1.63%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "157.89"
$ws.Range("E45").Value = "  +3.30%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.95"
$ws.Range("E46").Value = "  +4.82%  "

# Row 47
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.43"
$ws.Range("E47").Value = "  +2.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.04"
$ws.Range("E48").Value = "  +0.21%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.296"
$ws.Range("E49").Value = "  -0.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.71"
$ws.Range("E50").Value = "  -0.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.38"
$ws.Range("E51").Value = "  +0.02%  "
